$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the missing date value in row 4 (same date as row 3). Use Copy so the
# cell keeps referencing the same shared string / general formatting as A3
# instead of being re-interpreted (e.g. as a date serial number).
$ws.Range("A3").Copy($ws.Range("A4"))

# Update the long failure message shared string to the simpler "failed" text.
$ws.Range("D3").Value = "failed"
$ws.Range("D4").Value = "failed"

# Apply a red fill to the two "failed" result cells.
$ws.Range("D3").Interior.Color = 255
$ws.Range("D4").Interior.Color = 255
